# [Word] (TableCell) Map existing sample
# Append a new row to the "Snippets" table for the TableCell class sample.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

$newRow = 70

# --- Values -----------------------------------------------------------
$ws.Cells.Item($newRow, 1).Value = "TableCell"
# Column B (Method/Prop/Rel Name) is intentionally left blank for this row.
$ws.Cells.Item($newRow, 3).Value = "class"
$ws.Cells.Item($newRow, 4).Value = "word-tables-table-cell-access"
$ws.Cells.Item($newRow, 5).Value = "getTableCell"

# --- Formatting ---------------------------------------------------------
# Columns D/E reuse the same "vertical center" style already used by the
# rest of the data rows (e.g. row 69).
$ws.Range("D69").Copy()
$ws.Range("D70").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E69").Copy()
$ws.Range("E70").PasteSpecial(-4122)   # xlPasteFormats

# Columns A/B pick up the "General" number format explicitly applied to the
# new row header cells. Build that style on a scratch cell, copy it across,
# then clean the scratch cell back up.
$ws.Range("Z1").NumberFormat = "General"
$ws.Range("Z1").Copy()
$ws.Range("A70").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B70").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Z1").Clear()

# --- Table / list object -------------------------------------------------
$table = $ws.ListObjects.Item("Snippets")
$table.Resize($ws.Range("A1:E" + $newRow))

# --- View state: scroll so row 42 sits at the top of the frozen pane, with
# the newly added cell selected, matching the post-edit sheet view. -------
$ws.Activate()
$ws.Range("A42").Select()
$excel.ActiveWindow.ScrollRow = 42
$ws.Range("E" + $newRow).Select()
